# Rename the workbook's only sheet from "Property1" to "DataNode"
# ("unify the conception of DataNode, DataTable, Entity.")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"
